$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.662.88'
$ws.Range("E2").Value = '  +5.76%  '
$ws.Range("D3").Value = '3.064.37'
$ws.Range("E3").Value = '  +4.16%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.34'
$ws.Range("E5").Value = '  +5.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.18'
$ws.Range("E6").Value = '  +8.26%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.055.96'
$ws.Range("E8").Value = '  +3.94%  '
$ws.Range("E9").Value = '  +5.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  +3.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.22'
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("E12").Value = '  +5.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  +6.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.92'
$ws.Range("E14").Value = '  +6.90%  '
$ws.Range("D15").Value = '3.567.92'
$ws.Range("E15").Value = '  +4.00%  '
$ws.Range("D16").Value = '63.689.76'
$ws.Range("E16").Value = '  +5.61%  '
$ws.Range("D17").Value = '3.063.63'
$ws.Range("E17").Value = '  +3.98%  '
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.73'
$ws.Range("E19").Value = '  +5.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '484.05'
$ws.Range("E20").Value = '  +7.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("E21").Value = '  +5.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.685'
$ws.Range("E22").Value = '  +3.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.20'
$ws.Range("E23").Value = '  +7.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.70'
$ws.Range("E24").Value = '  +5.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.57'
$ws.Range("E25").Value = '  +9.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("E27").Value = '  +6.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.96'
$ws.Range("E28").Value = '  +5.49%  '
$ws.Range("E29").Value = '  +10.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.07'
$ws.Range("E31").Value = '  +4.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  +3.24%  '
$ws.Range("E33").Value = '  +9.48%  '
$ws.Range("E34").Value = '  +10.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.39'
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.98'
$ws.Range("E36").Value = '  +5.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '468.11'
$ws.Range("E37").Value = '  +5.56%  '
$ws.Range("D38").Value = '3.160.08'
$ws.Range("E38").Value = '  +1.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0813'
$ws.Range("E39").Value = '  +6.56%  '
$ws.Range("E40").Value = '  +5.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.120'
$ws.Range("E41").Value = '  +3.73%  '
$ws.Range("E42").Value = '  +4.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '28.46'
$ws.Range("E43").Value = '  +15.41%  '
$ws.Range("E44").Value = '  +6.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.252'
$ws.Range("E45").Value = '  +4.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.04'
$ws.Range("E47").Value = '  +7.87%  '
$ws.Range("E48").Value = '  +3.66%  '
$ws.Range("D49").Value = '0.0₃0512'
$ws.Range("E49").Value = '  +3.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '116.07'
$ws.Range("E50").Value = '  -1.75%  '
$ws.Range("E51").Value = '  +6.90%  '
